$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the per-row IP addresses in column F (rows 3-6) with the same
# localhost address already used in F2, collapsing the five distinct
# "192.168.1.11x" shared strings down to a single "127.0.0.1" entry.
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("F3").Value = "127.0.0.1"
$ws.Range("F4").Value = "127.0.0.1"
$ws.Range("F5").Value = "127.0.0.1"
$ws.Range("F6").Value = "127.0.0.1"

# Update the saved selection/active cell on the sheet view.
$ws.Range("F14").Select()
